# Wave 5 React Tracker - add two new progress rows (26/08/2016 and 27/08/2016)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: date, Compose Mail = Completed, View Mail = In progress
$ws.Range("A7").Value = "26/08/2016"
$ws.Range("I7").Value = "Completed"
$ws.Range("J7").Value = "In progress"

# Row 8: date, View Mail = Completed
$ws.Range("A8").Value = "27/08/2016"
$ws.Range("J8").Value = "Completed"

# Match the updated view/selection state
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("M8").Select()
